# Insert a new data row at row 175 ("Hortaliza, Macroferia Regional de Talca - Zapallo italiano").
# This pushes the former rows 175-288 down to 176-289, and the new row 175
# receives a fresh record (Fecha 2022-02-11 / serial 44603).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A175").EntireRow.Insert()

$ws.Cells.Item(175, 1).Value  = 5
$ws.Cells.Item(175, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(175, 3).Value  = "Maule"
$ws.Cells.Item(175, 4).Value  = 44603
$ws.Cells.Item(175, 5).Value  = 7
$ws.Cells.Item(175, 6).Value  = 100112032
$ws.Cells.Item(175, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(175, 8).Value  = "Sin especificar"
$ws.Cells.Item(175, 9).Value  = "Primera"
$ws.Cells.Item(175, 10).Value = 400
$ws.Cells.Item(175, 11).Value = 5000
$ws.Cells.Item(175, 12).Value = 5000
$ws.Cells.Item(175, 13).Value = 5000
$ws.Cells.Item(175, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(175, 15).Value = "Región del Maule"
$ws.Cells.Item(175, 16).Value = 100
$ws.Cells.Item(175, 17).Value = 50
$ws.Cells.Item(175, 18).Value = "Hortaliza"
